# FINFLUX-2698 Correcting Overdue failed scenarios
#
# This script replays, via Excel COM-interop, the "overdue scenario
# correction" edit captured in the target diff:
#   - Summary sheet: small rounding correction (195.49 -> 195.53)
#   - Repayment schedule: several recalculated installment values
#   - Transactions sheet: two additional transactions (accrual / income
#     posting) were captured, shifting the disbursement row further down
#     and renumbering every transaction id
#   - The active worksheet moves from "ChargesTab" to "Transactions"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: rounding correction
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A3").Value = 195.53
$wsSummary.Range("E3").Value = 195.53

# ---------------------------------------------------------------------
# Repayment schedule: corrected installment figures
# ---------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Range("F4").Value = 841.22
$wsRepay.Range("G4").Value = 3322.07
$wsRepay.Range("H4").Value = 46.5

$wsRepay.Range("F5").Value = 838.08
$wsRepay.Range("G5").Value = 2483.9899999999998
$wsRepay.Range("H5").Value = 49.64

$wsRepay.Range("F6").Value = 863.23
$wsRepay.Range("G6").Value = 1620.76
$wsRepay.Range("H6").Value = 24.49

$wsRepay.Range("G7").Value = 749.55

$wsRepay.Range("F8").Value = 749.55
$wsRepay.Range("K8").Value = 756.93
$wsRepay.Range("Q8").Value = 756.93

# ---------------------------------------------------------------------
# Transactions sheet: two newly captured transactions push everything
# down; ids / dates / amounts / running balances all get renumbered.
# ---------------------------------------------------------------------
$wsTx = $wb.Worksheets.Item("Transactions")

# Insert two rows above the disbursement row (old row 8) so it lands on
# row 10 - Excel clones the formatting of the row above for the newly
# inserted rows, and keeps the disbursement row's own formatting intact.
$wsTx.Rows("8:9").Insert()

# Row 2 (Income Posting)
$wsTx.Range("A2").Value = 203
$wsTx.Range("E2").Value = 51.91
$wsTx.Range("J2").Value = 5145.3999999999996

# Row 3 (Accrual)
$wsTx.Range("A3").Value = 202
$wsTx.Range("E3").Value = 0.03

# Row 4 (Accrual) - also moved to the later (42094) accrual date
$wsTx.Range("A4").Value = 198
$wsTx.Range("C4").Value = 42094
$wsTx.Range("E4").Value = 51.88

# Row 5 (Income Posting)
$wsTx.Range("A5").Value = 201
$wsTx.Range("E5").Value = 49.11
$wsTx.Range("J5").Value = 5093.49

# Row 6 - becomes an Accrual row (was Income Posting); copy the format
# from row 7 (an existing Accrual row) before setting the value so the
# balance column switches from currency back to the plain/general style.
$wsTx.Range("A6").Value = 200
$wsTx.Range("C6").Value = 42063
$wsTx.Range("D6").Value = "Accrual"
$wsTx.Range("E6").Value = 2.68
$wsTx.Range("J7").Copy()
$wsTx.Range("J6").PasteSpecial(-4122)
$wsTx.Range("J6").Value = 0

# Row 7 (Accrual)
$wsTx.Range("A7").Value = 196
$wsTx.Range("C7").Value = 42063
$wsTx.Range("E7").Value = 46.43

# Row 8 (newly inserted row) - becomes an Income Posting row; copy the
# currency format from row 5 (an existing Income Posting row) first.
$wsTx.Range("A8").Value = 195
$wsTx.Range("B8").Value = "Head Office"
$wsTx.Range("C8").Value = 42035
$wsTx.Range("D8").Value = "Income Posting"
$wsTx.Range("E8").Value = 44.38
$wsTx.Range("J5").Copy()
$wsTx.Range("J8").PasteSpecial(-4122)
$wsTx.Range("J8").Value = 5044.38

# Row 9 (newly inserted row) - Accrual row, already has the right
# (general/plain) formatting inherited from the Insert() above.
$wsTx.Range("A9").Value = 194
$wsTx.Range("B9").Value = "Head Office"
$wsTx.Range("C9").Value = 42035
$wsTx.Range("D9").Value = "Accrual"
$wsTx.Range("E9").Value = 44.38
$wsTx.Range("J9").Value = 0

# Row 10 - this is the original disbursement row, shifted down by the
# insert above with all of its values/formatting intact; only its id
# changes.
$wsTx.Range("A10").Value = 193

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Selections: restore the cursor position on each sheet the way the
# author left them, and finish with "Transactions" as the active tab
# (it was "ChargesTab" before).
# ---------------------------------------------------------------------
$wsSummary.Range("C10").Select()
$wsRepay.Range("L5").Select()
$wb.Worksheets.Item("ChargesTab").Range("D6").Select()
$wsTx.Range("H4").Select()
